$d = $word.ActiveDocument

# Step 1: update the text of the first paragraph's run, adding two trailing spaces.
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a Microsoft word document.  ", 2)

# Step 2: insert a new run right after the updated text (but before the
# paragraph mark), colored dark red, containing the change-notice text.
$para = $d.Paragraphs(1)
$insertRange = $para.Range
[void]$insertRange.MoveEnd(1, -1)   # wdCharacter, -1: exclude the paragraph mark
$insertRange.Collapse(0)      # wdCollapseEnd
$insertRange.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch alternate)")
$insertRange.Font.Color = 192  # RGB(192,0,0) -> 0x0000C0 as Word's BGR-packed long
